$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2171428571428571
$ws.Range("C2").Value = 0.5171428571428571
$ws.Range("J2").Value = 0.01714285714285714
$ws.Range("P2").Value = 0.1685714285714286
$ws.Range("S2").Value = 0.08
$ws.Range("B3").Value = 0.01630434782608696
$ws.Range("C3").Value = 0.005434782608695652
$ws.Range("J3").Value = 0.01630434782608696
$ws.Range("P3").Value = 0.7663043478260869
$ws.Range("S3").Value = 0.1956521739130435
$ws.Range("J4").Value = 0.01818181818181818
$ws.Range("P4").Value = 0.7454545454545455
$ws.Range("S4").Value = 0.2363636363636364
$ws.Range("B6").Value = 0.09392265193370165
$ws.Range("D6").Value = 0.02209944751381215
$ws.Range("F6").Value = 0.04972375690607735
$ws.Range("J6").Value = 0.2541436464088398
$ws.Range("O6").Value = 0.01104972375690608
$ws.Range("Q6").Value = 0.1546961325966851
$ws.Range("R6").Value = 0.07734806629834254
$ws.Range("S6").Value = 0.3370165745856354
$ws.Range("B7").Value = 0.0776255707762557
$ws.Range("D7").Value = 0.0273972602739726
$ws.Range("F7").Value = 0.045662100456621
$ws.Range("J7").Value = 0.1415525114155251
$ws.Range("O7").Value = 0.0045662100456621
$ws.Range("Q7").Value = 0.1598173515981735
$ws.Range("R7").Value = 0.0821917808219178
$ws.Range("S7").Value = 0.4611872146118721
$ws.Range("B8").Value = 0.1053811659192825
$ws.Range("D8").Value = 0.0179372197309417
$ws.Range("F8").Value = 0.04708520179372197
$ws.Range("J8").Value = 0.1412556053811659
$ws.Range("O8").Value = 0.01569506726457399
$ws.Range("Q8").Value = 0.1614349775784753
$ws.Range("R8").Value = 0.09417040358744394
$ws.Range("S8").Value = 0.4170403587443946
$ws.Range("B9").Value = 0.1157407407407407
$ws.Range("D9").Value = 0.02314814814814815
$ws.Range("F9").Value = 0.05092592592592592
$ws.Range("J9").Value = 0.1018518518518518
$ws.Range("O9").Value = 0.02777777777777778
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.1018518518518518
$ws.Range("S9").Value = 0.412037037037037
$ws.Range("B10").Value = 0.1334367726920093
$ws.Range("D10").Value = 0.02560124127230411
$ws.Range("F10").Value = 0.05663304887509697
$ws.Range("J10").Value = 0.108611326609775
$ws.Range("O10").Value = 0.01551590380139643
$ws.Range("Q10").Value = 0.2048099301784329
$ws.Range("R10").Value = 0.09154383242823895
$ws.Range("S10").Value = 0.3638479441427463
$ws.Range("G11").Value = 0.1221122112211221
$ws.Range("J11").Value = 0.0924092409240924
$ws.Range("K11").Value = 0.1782178217821782
$ws.Range("L11").Value = 0.5907590759075908
$ws.Range("S11").Value = 0.0165016501650165
$ws.Range("G12").Value = 0.7788944723618091
$ws.Range("J12").Value = 0.1507537688442211
$ws.Range("K12").Value = 0.005025125628140704
$ws.Range("L12").Value = 0.03517587939698492
$ws.Range("S12").Value = 0.03015075376884422
$ws.Range("G13").Value = 0.7272727272727273
$ws.Range("J13").Value = 0.1590909090909091
$ws.Range("S13").Value = 0.1136363636363636
$ws.Range("F15").Value = 0.01818181818181818
$ws.Range("H15").Value = 0.1272727272727273
$ws.Range("I15").Value = 0.09090909090909091
$ws.Range("J15").Value = 0.4090909090909091
$ws.Range("K15").Value = 0.05909090909090909
$ws.Range("M15").Value = 0.004545454545454545
$ws.Range("O15").Value = 0.07727272727272727
$ws.Range("S15").Value = 0.2136363636363636
$ws.Range("F16").Value = 0.01265822784810127
$ws.Range("H16").Value = 0.1181434599156118
$ws.Range("I16").Value = 0.0759493670886076
$ws.Range("J16").Value = 0.4978902953586498
$ws.Range("K16").Value = 0.06329113924050633
$ws.Range("M16").Value = 0.01687763713080169
$ws.Range("N16").Value = 0.004219409282700422
$ws.Range("O16").Value = 0.0970464135021097
$ws.Range("S16").Value = 0.1139240506329114
$ws.Range("F17").Value = 0.004608294930875576
$ws.Range("H17").Value = 0.1658986175115207
$ws.Range("I17").Value = 0.08294930875576037
$ws.Range("J17").Value = 0.4331797235023042
$ws.Range("K17").Value = 0.1013824884792627
$ws.Range("M17").Value = 0.02764976958525346
$ws.Range("O17").Value = 0.05529953917050692
$ws.Range("S17").Value = 0.1290322580645161
$ws.Range("F18").Value = 0.01408450704225352
$ws.Range("H18").Value = 0.1784037558685446
$ws.Range("I18").Value = 0.0892018779342723
$ws.Range("J18").Value = 0.3615023474178404
$ws.Range("K18").Value = 0.0892018779342723
$ws.Range("M18").Value = 0.0187793427230047
$ws.Range("N18").Value = 0.004694835680751174
$ws.Range("O18").Value = 0.09389671361502347
$ws.Range("S18").Value = 0.1502347417840376
$ws.Range("F19").Value = 0.01879699248120301
$ws.Range("H19").Value = 0.2022556390977444
$ws.Range("I19").Value = 0.09323308270676692
$ws.Range("J19").Value = 0.3548872180451128
$ws.Range("K19").Value = 0.1165413533834586
$ws.Range("M19").Value = 0.01729323308270677
$ws.Range("O19").Value = 0.05413533834586466
$ws.Range("S19").Value = 0.1428571428571428
